$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format ("@") on Price cells whose new value would otherwise
# be auto-parsed by Excel as a number, so they stay text like the source data.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "63.165.73"
$ws.Range("E2").Value = "  -0.80%  "

# Row 3
$ws.Range("D3").Value = "3.186.93"
$ws.Range("E3").Value = "  -3.80%  "

# Row 4
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").Value = "592.64"
$ws.Range("E5").Value = "  -1.97%  "

# Row 6
$ws.Range("D6").Value = "135.53"
$ws.Range("E6").Value = "  -4.04%  "

# Row 7
$ws.Range("E7").Value = "  -0.05%  "

# Row 8
$ws.Range("D8").Value = "3.185.33"
$ws.Range("E8").Value = "  -3.86%  "

# Row 9
$ws.Range("E9").Value = "  -0.57%  "

# Row 10
$ws.Range("D10").Value = "0.142"
$ws.Range("E10").Value = "  -5.57%  "

# Row 11
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -4.93%  "

# Row 12
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").Value = "  -2.77%  "

# Row 13
$ws.Range("D13").Value = "0.0000238"
$ws.Range("E13").Value = "  -4.10%  "

# Row 14
$ws.Range("D14").Value = "34.81"
$ws.Range("E14").Value = "  -0.34%  "

# Row 15
$ws.Range("D15").Value = "3.717.99"
$ws.Range("E15").Value = "  -3.67%  "

# Row 16
$ws.Range("E16").Value = "  -0.99%  "

# Row 17
$ws.Range("D17").Value = "3.185.86"
$ws.Range("E17").Value = "  -3.93%  "

# Row 18
$ws.Range("D18").Value = "63.138.86"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19
$ws.Range("D19").Value = "6.59"
$ws.Range("E19").Value = "  -3.82%  "

# Row 20
$ws.Range("D20").Value = "463.27"
$ws.Range("E20").Value = "  -3.58%  "

# Row 21
$ws.Range("D21").Value = "14.05"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  -5.45%  "

# Row 23
$ws.Range("D23").Value = "7.66"
$ws.Range("E23").Value = "  -4.20%  "

# Row 24
$ws.Range("D24").Value = "13.43"
$ws.Range("E24").Value = "  -4.26%  "

# Row 25
$ws.Range("D25").Value = "82.95"
$ws.Range("E25").Value = "  -2.72%  "

# Row 26
$ws.Range("E26").Value = "  +0.11%  "

# Row 27
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "2.68"
$ws.Range("E27").Value = "  -3.43%  "

# Row 28
$ws.Range("B28").Value = "FirstDigitalUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.03%  "

# Row 29
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "7.74"
$ws.Range("E29").Value = "  -5.73%  "

# Row 30
$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "6.80"
$ws.Range("E30").Value = "  -5.34%  "

# Row 31
$ws.Range("D31").Value = "2.04"
$ws.Range("E31").Value = "  -4.83%  "

# Row 32
$ws.Range("D32").Value = "27.30"
$ws.Range("E32").Value = "  -5.50%  "

# Row 33
$ws.Range("E33").Value = "  -2.90%  "

# Row 34
$ws.Range("D34").Value = "2.38"
$ws.Range("E34").Value = "  -5.28%  "

# Row 35
$ws.Range("D35").Value = "1.04"
$ws.Range("E35").Value = "  -6.22%  "

# Row 36
$ws.Range("D36").Value = "5.84"
$ws.Range("E36").Value = "  -3.85%  "

# Row 37
$ws.Range("D37").Value = "51.35"
$ws.Range("E37").Value = "  -2.05%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0710"
$ws.Range("E38").Value = "  -4.57%  "

# Row 39
$ws.Range("D39").Value = "0.0390"
$ws.Range("E39").Value = "  -2.51%  "

# Row 40
$ws.Range("D40").Value = "407.69"
$ws.Range("E40").Value = "  -5.82%  "

# Row 41
$ws.Range("D41").Value = "8.11"
$ws.Range("E41").Value = "  -2.76%  "

# Row 42
$ws.Range("D42").Value = "2.66"
$ws.Range("E42").Value = "  -2.64%  "

# Row 43
$ws.Range("D43").Value = "0.113"
$ws.Range("E43").Value = "  -5.15%  "

# Row 44
$ws.Range("D44").Value = "2.808.42"
$ws.Range("E44").Value = "  -10.11%  "

# Row 45
$ws.Range("D45").Value = "0.253"
$ws.Range("E45").Value = "  -5.00%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.13"
$ws.Range("E46").Value = "  -4.75%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.02%  "

# Row 48
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").Value = "35.31"
$ws.Range("E48").Value = "  -3.78%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "124.45"
$ws.Range("E49").Value = "  +0.12%  "

# Row 50
$ws.Range("D50").Value = "25.33"
$ws.Range("E50").Value = "  -3.72%  "

# Row 51
$ws.Range("D51").Value = "0.112"
$ws.Range("E51").Value = "  -1.59%  "
